$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F11").Value = "Account Number = 709`nClient Number = 9710`nBalance = 450`nDate Created = July 9, 2018`nmanagement_fee = 2"
$ws.Range("F13").Value = "Account Number = 709`nClient Number = 9710`nBalance = 450`nDate Created = July 9, 2018`nmanagement_fee = 2"
$ws.Range("F10").Value = "Account Number = 709`nClient Number = 9710`nBalance = 450`nDate Created = July 9, 2014`nmanagement_fee = 2"
$ws.Range("F7").Value = "Account Number = 709`nClient Number = 9710`nBalance = 450`nDate Created = July 9, 2011`nmanagement_fee = 2"
$ws.Range("F9").Value = "Account Number = 709`nClient Number = 9710`nBalance = 450`nDate Created = July 9, 2011`nmanagement_fee = 2"
$ws.Range("F12").Value = "Account Number = 709`nClient Number = 9710`nBalance = 450`nDate Created = July 9, 2011`nmanagement_fee = 2"
